# Auto-generated edit script
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the numeric/percent-looking strings in columns D and E stay as
# plain text (matching the inlineStr cells already in the sheet) instead
# of being auto-converted to numbers by Excel.
$ws.Range("D2:E50").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '291.76'
$ws.Range("E2").Value = '-0.41%'

# Row 3
$ws.Range("D3").Value = '40.25'
$ws.Range("E3").Value = '0.49%'

# Row 4
$ws.Range("D4").Value = '5.006'
$ws.Range("E4").Value = '-0.79%'

# Row 5
$ws.Range("D5").Value = '0.07290'
$ws.Range("E5").Value = '-1.61%'

# Row 6
$ws.Range("B6").Value = 'GateToken'
$ws.Range("C6").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D6").Value = '4.284'
$ws.Range("E6").Value = '-1.41%'

# Row 7
$ws.Range("B7").Value = 'FTXToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D7").Value = '1.551'
$ws.Range("E7").Value = '-1.60%'

# Row 8
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").Value = '0.9254'
$ws.Range("E8").Value = '-0.59%'

# Row 9
$ws.Range("B9").Value = 'BTSEToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D9").Value = '2.361'
$ws.Range("E9").Value = '-2.44%'

# Row 10
$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D10").Value = '0.1157'
$ws.Range("E10").Value = '-2.56%'

# Row 11
$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D11").Value = '0.1760'
$ws.Range("E11").Value = '1.09%'

# Row 12
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").Value = '0.04351'
$ws.Range("E12").Value = '3.62%'

# Row 13
$ws.Range("B13").Value = 'MandalaExchangeToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D13").Value = '0.08686'
$ws.Range("E13").Value = '-0.56%'

# Row 14
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").Value = '0.1053'
$ws.Range("E14").Value = '-0.18%'

# Row 15
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").Value = '0.001278'
$ws.Range("E15").Value = '0.07%'

# Row 16
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").Value = '0.005914'
$ws.Range("E16").Value = '-0.20%'

# Row 17
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").Value = '3.341'
$ws.Range("E17").Value = '-0.62%'

# Row 18
$ws.Range("D18").Value = '0.3284'
$ws.Range("E18").Value = '-1.96%'

# Row 19
$ws.Range("D19").Value = '7.825'
$ws.Range("E19").Value = '1.54%'

# Row 20
$ws.Range("D20").Value = '0.1391'
$ws.Range("E20").Value = '1.92%'

# Row 21
$ws.Range("E21").Value = '-1.73%'

# Row 22
$ws.Range("D22").Value = '0.03921'
$ws.Range("E22").Value = '1.24%'

# Row 23
$ws.Range("E23").Value = '-2.83%'

# Row 24
$ws.Range("D24").Value = '0.003647'
$ws.Range("E24").Value = '3.51%'

# Row 25
$ws.Range("D25").Value = '0.0001202'
$ws.Range("E25").Value = '-8.19%'

# Row 26
$ws.Range("D26").Value = '0.0003727'
$ws.Range("E26").Value = '-0.94%'

# Row 38
$ws.Range("D38").Value = '0.02308'
$ws.Range("E38").Value = '0.23%'

# Row 39
$ws.Range("D39").Value = '0.05065'
$ws.Range("E39").Value = '1.43%'

# Row 40
$ws.Range("D40").Value = '0.005748'
$ws.Range("E40").Value = '37.87%'

# Row 41
$ws.Range("D41").Value = '0.007857'
$ws.Range("E41").Value = '2.07%'

# Row 42
$ws.Range("D42").Value = '0.1286'
$ws.Range("E42").Value = '0.41%'

# Row 43
$ws.Range("D43").Value = '0.007396'
$ws.Range("E43").Value = '-2.29%'

# Row 44
$ws.Range("D44").Value = '0.007265'
$ws.Range("E44").Value = '0.55%'

# Row 45
$ws.Range("D45").Value = '0.3183'
$ws.Range("E45").Value = '0.43%'

# Row 46
$ws.Range("D46").Value = '0.00006206'
$ws.Range("E46").Value = '-6.20%'

# Row 47
$ws.Range("E47").Value = '-0.94%'

# Row 48
$ws.Range("D48").Value = '0.04835'
$ws.Range("E48").Value = '-80.80%'

# Row 49
$ws.Range("D49").Value = '0.00002103'
$ws.Range("E49").Value = '-0.94%'

# Row 50
$ws.Range("D50").Value = '0.0002003'
$ws.Range("E50").Value = '-0.94%'
